# This script re-applies the row reordering that the commit performed on
# "纺织服装、服饰业工业生产者出厂价格指数(上年同月=100).xlsx".
#
# Within every calendar year block the rows for Oct/Nov/Dec were moved to
# the top of that year's block (ahead of Jan-Sep), while the (date, value)
# pairs themselves were left untouched - only their row position changed.
# The table below lists, for every data row (2-49), the final row number
# together with the (Date, B, C, D) values that must end up there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  Date="2014-10"; B=100.8146;  C=99.9603;   D=$null },
    @{ Row=3;  Date="2014-11"; B=101.0944;  C=100.0519;  D=$null },
    @{ Row=4;  Date="2014-12"; B=101.3665;  C=100.1502;  D=$null },
    @{ Row=5;  Date="2014-01"; B=102.1243;  C=100.2447;  D=$null },
    @{ Row=6;  Date="2014-02"; B=101.2312;  C=100.3575;  D=$null },
    @{ Row=7;  Date="2014-03"; B=101.4826;  C=100.268;   D=$null },
    @{ Row=8;  Date="2014-04"; B=101.8418;  C=100.2129;  D=$null },
    @{ Row=9;  Date="2014-05"; B=101.6706;  C=100.2494;  D=$null },
    @{ Row=10; Date="2014-06"; B=102.2387;  C=100.2766;  D=$null },
    @{ Row=11; Date="2014-07"; B=102.0669;  C=100.2108;  D=$null },
    @{ Row=12; Date="2014-08"; B=101.4958;  C=100.2296;  D=$null },
    @{ Row=13; Date="2014-09"; B=101.5385;  C=100.0181;  D=$null },
    @{ Row=14; Date="2015-10"; B=102.4;     C=100.9;     D=$null },
    @{ Row=15; Date="2015-11"; B=102.624;   C=100.7543;  D=$null },
    @{ Row=16; Date="2015-12"; B=102.9569;  C=100.5755;  D=$null },
    @{ Row=17; Date="2015-01"; B=101.001;   C=100.6496;  D=$null },
    @{ Row=18; Date="2015-02"; B=102.3155;  C=100.6173;  D=$null },
    @{ Row=19; Date="2015-03"; B=101.8702;  C=100.8709;  D=$null },
    @{ Row=20; Date="2015-04"; B=101.7441;  C=100.6292;  D=$null },
    @{ Row=21; Date="2015-05"; B=101.1222;  C=100.5848;  D=$null },
    @{ Row=22; Date="2015-06"; B=102.1945;  C=100.6486;  D=$null },
    @{ Row=23; Date="2015-07"; B=101.7685;  C=100.7167;  D=$null },
    @{ Row=24; Date="2015-08"; B=102.06;    C=100.7071;  D=$null },
    @{ Row=25; Date="2015-09"; B=102.4824;  C=100.8275;  D=$null },
    @{ Row=26; Date="2016-10"; B=100.6;     C=100.1;     D=101 },
    @{ Row=27; Date="2016-11"; B=100.6;     C=100.4;     D=101.3 },
    @{ Row=28; Date="2016-12"; B=100.8;     C=100.9;     D=101.8 },
    @{ Row=29; Date="2016-01"; B=102.6548;  C=100.3357;  D=100.5647 },
    @{ Row=30; Date="2016-02"; B=101.696;   C=100.253;   D=100.1751 },
    @{ Row=31; Date="2016-03"; B=101.9562;  C=100.1799;  D=100.6021 },
    @{ Row=32; Date="2016-04"; B=101.6713;  C=100.3656;  D=100.6055 },
    @{ Row=33; Date="2016-05"; B=101.7;     C=100.4;     D=100.4 },
    @{ Row=34; Date="2016-06"; B=100.9;     C=100.5;     D=100.9 },
    @{ Row=35; Date="2016-07"; B=101.4;     C=100.5;     D=101.1 },
    @{ Row=36; Date="2016-08"; B=100.7;     C=100.3;     D=100.9 },
    @{ Row=37; Date="2016-09"; B=100.3;     C=100.2;     D=100.8 },
    @{ Row=38; Date="2017-10"; B=100.4;     C=100.8;     D=100.4 },
    @{ Row=39; Date="2017-11"; B=100.5;     C=100.8;     D=100.3 },
    @{ Row=40; Date="2017-12"; B=100.2;     C=100.6;     D=99.9 },
    @{ Row=41; Date="2017-01"; B=101;       C=100.5;     D=101.1 },
    @{ Row=42; Date="2017-02"; B=101;       C=100.8;     D=101.6 },
    @{ Row=43; Date="2017-03"; B=100.9;     C=100.8;     D=101.7 },
    @{ Row=44; Date="2017-04"; B=101.2;     C=100.9;     D=101.5 },
    @{ Row=45; Date="2017-05"; B=100.9;     C=100.8;     D=101.6 },
    @{ Row=46; Date="2017-06"; B=100.9;     C=100.8;     D=101.2 },
    @{ Row=47; Date="2017-07"; B=100.3;     C=100.8;     D=100.9 },
    @{ Row=48; Date="2017-08"; B=100.4;     C=101.1;     D=100.9 },
    @{ Row=49; Date="2017-09"; B=100.2;     C=100.9;     D=100.6 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($null -ne $r.D) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
}
